# Daily attendance processing - 2026-01-18 14:59:09
# Swap the order of "Recorded By" entries in column G where the value is a
# two-part, comma-separated list starting with "dnasr281@gmail.com" so that
# it appears last instead of first (e.g. "dnasr281@gmail.com, System" ->
# "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $value = $cell.Value()

    if ($value -ne $null -and $value -like "dnasr281@gmail.com,*") {
        $parts = $value -split ", ", 2
        if ($parts.Count -eq 2) {
            $newValue = $parts[1] + ", " + $parts[0]
            $cell.Value = $newValue
        }
    }
}
